# Insert a new data row at row 447 (pushing all existing rows 447..506 down to
# 448..507) and populate it with the new record. Everything below shifts down
# automatically via the native Insert, preserving formatting (e.g. the date
# style on column D) and all other row contents/styles untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(447).Insert()

$ws.Range("A447").Value = 3
$ws.Range("B447").Value = "Femacal de La Calera"
$ws.Range("C447").Value = "Coquimbo"
$ws.Range("D447").Value = "2023-02-27"
$ws.Range("E447").Value = 5
$ws.Range("F447").Value = 100114013
$ws.Range("G447").Value = "Zanahoria"
$ws.Range("H447").Value = "Sin especificar"
$ws.Range("I447").Value = "Primera"
$ws.Range("J447").Value = 510
$ws.Range("K447").Value = 9000
$ws.Range("L447").Value = 9500
$ws.Range("M447").Value = 9225
$ws.Range("N447").Value = "$/saco 20 kilos"
$ws.Range("O447").Value = "Provincia de Quillota"
$ws.Range("P447").Value = 461
$ws.Range("Q447").Value = 20
$ws.Range("R447").Value = "Hortaliza"
